# Apply "Atualizacao de bases das ligas" update to the Lithuania A Lyga sheet.
# The underlying change re-orders several match rows (the row index in column A
# stays fixed, but all the other match data in columns B:AC moves to a
# different row). This script reads each affected row's B:AC values first,
# then writes them back in the new order, so that no data is lost while the
# rows are being rearranged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues([int]$row) {
    return $ws.Range("B" + $row + ":AC" + $row).Value2
}

function Set-RowValues([int]$row, $values) {
    $ws.Range("B" + $row + ":AC" + $row).Value2 = $values
}

# --- Simple two-row swaps -------------------------------------------------
$swapPairs = @(
    @(26, 27),
    @(50, 51),
    @(117, 118),
    @(136, 137)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# --- Five-row cycle (rows 100-104) ----------------------------------------
# New content per row comes from the following previous rows:
#   row 100 <- old row 103
#   row 101 <- old row 102
#   row 102 <- old row 100
#   row 103 <- old row 104
#   row 104 <- old row 101
$v100 = Get-RowValues 100
$v101 = Get-RowValues 101
$v102 = Get-RowValues 102
$v103 = Get-RowValues 103
$v104 = Get-RowValues 104

Set-RowValues 100 $v103
Set-RowValues 101 $v102
Set-RowValues 102 $v100
Set-RowValues 103 $v104
Set-RowValues 104 $v101
